$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the "Name" property value (row 4, column B) that was previously empty.
$ws.Range("B4").Value = "LieunaissanceVs"

# Update the "Date" property value (row 8, column B) to reflect the regeneration timestamp.
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
